$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.96404131162324602
$ws.Range("C1").Value = 0.9722774741644733
$ws.Range("C2").Value = 0.69716924076987286
$ws.Range("E3").Value = 0.91162545894547375
$ws.Range("C4").Value = 0.83611591316431222
$ws.Range("W4").Value = 0.98215011412982878
$ws.Range("D5").Value = 0.99468642265767759
$ws.Range("F5").Value = 0.73607571893886981
$ws.Range("D6").Value = 0.71955889109527915
$ws.Range("G6").Value = 0.97757198793395816
$ws.Range("AP6").Value = 0.74442793924016737
$ws.Range("H7").Value = 0.9831269269812013
$ws.Range("I7").Value = 0.81009540416698655
$ws.Range("K7").Value = 0.65979517431239088
$ws.Range("J8").Value = 0.89068915732884313
$ws.Range("H9").Value = 0.7754406674236034
$ws.Range("K10").Value = 0.9176088613743929
$ws.Range("I11").Value = 0.6574104389980111
$ws.Range("M11").Value = 0.76399938584775584
$ws.Range("AL11").Value = 0.9684928643203633
$ws.Range("A12").Value = 0.65434195644310766
$ws.Range("N12").Value = 0.60393266527644662
$ws.Range("B13").Value = 0.92304419887958011
$ws.Range("O13").Value = 0.86058445079850365
$ws.Range("A14").Value = 0.93266364949188207
$ws.Range("P14").Value = 0.90545177207171257
$ws.Range("BG15").Value = 0.94505277871758719
$ws.Range("AV16").Value = 0.7942107292693773
$ws.Range("O17").Value = 0.95887040675724799
$ws.Range("R17").Value = 0.95174550026614135
$ws.Range("S17").Value = 0.82946641892017836
$ws.Range("P18").Value = 0.89026985570489636
$ws.Range("T18").Value = 0.83461791423801002
$ws.Range("BC18").Value = 0.88747840223790542
$ws.Range("R19").Value = 0.77848112039344619
$ws.Range("T19").Value = 0.77105852359532046
$ws.Range("AN19").Value = 0.95116226987349761
$ws.Range("V20").Value = 0.88388577050429218
$ws.Range("T21").Value = 0.69265382805122777
$ws.Range("AO21").Value = 0.69545849899330392
$ws.Range("U22").Value = 0.91040171847966711
$ws.Range("Y22").Value = 0.98131202708736276
$ws.Range("V23").Value = 0.91687330530020872
$ws.Range("X23").Value = 0.88909909907635964
$ws.Range("V24").Value = 0.9723767349245005
$ws.Range("Z25").Value = 0.88992835198492859
$ws.Range("X26").Value = 0.70387076374587498
$ws.Range("AX26").Value = 0.9897573031392547
$ws.Range("BB26").Value = 0.77562266731329077
$ws.Range("H27").Value = 0.93558944026620749
$ws.Range("AB27").Value = 0.81052698512537824
$ws.Range("BO27").Value = 0.95864385437994715
$ws.Range("Z28").Value = 0.59949099136698747
$ws.Range("AD28").Value = 0.8331732874395994
$ws.Range("AO28").Value = 0.64692846680566385
$ws.Range("AE29").Value = 0.88574204788402744
$ws.Range("AC30").Value = 0.92618773991791992
$ws.Range("AE30").Value = 0.94325677749818471
$ws.Range("AG31").Value = 0.85684059416981961
$ws.Range("AZ31").Value = 0.97510819433751017
$ws.Range("AD32").Value = 0.9554707839651877
$ws.Range("AG32").Value = 0.84890146257408838
$ws.Range("AH32").Value = 0.86045026925215695
$ws.Range("BD32").Value = 0.72193384048071363
$ws.Range("K33").Value = 0.77773856047530443
$ws.Range("AJ34").Value = 0.99516783717125135
$ws.Range("AG35").Value = 0.62468392182231725
$ws.Range("AH35").Value = 0.7666453645483513
$ws.Range("AI36").Value = 0.84962424822352622
$ws.Range("AK36").Value = 0.94750116522584149
$ws.Range("BP36").Value = 0.97368444105157226
$ws.Range("AI37").Value = 0.9996976293494706
$ws.Range("AM37").Value = 0.97621110182866633
$ws.Range("B38").Value = 0.95876438699931787
$ws.Range("X38").Value = 0.67690109308001245
$ws.Range("AJ38").Value = 0.90043444177007614
$ws.Range("AK38").Value = 0.93512821315573935
$ws.Range("AL39").Value = 0.7650368374631058
$ws.Range("AO39").Value = 0.94155887333972932
$ws.Range("AM40").Value = 0.87517773023156531
$ws.Range("AO40").Value = 0.76089656725935617
$ws.Range("AP40").Value = 0.62963159304249317
$ws.Range("AT41").Value = 0.85063640616951952
$ws.Range("AQ42").Value = 0.84943002081235552
$ws.Range("AO43").Value = 0.98945383120922725
$ws.Range("AS43").Value = 0.60887744139869571
$ws.Range("AU43").Value = 0.83686503628319098
$ws.Range("AP44").Value = 0.83650807412143924
$ws.Range("AS44").Value = 0.95359281320346112
$ws.Range("AT44").Value = 0.66402263764796321
$ws.Range("AT45").Value = 0.91941263108681248
$ws.Range("AU45").Value = 0.85745813629046985
$ws.Range("AW47").Value = 0.94778582011574164
$ws.Range("AT48").Value = 0.65718604753677434
$ws.Range("AU48").Value = 0.98318764357021593
$ws.Range("AW48").Value = 0.90823185693024044
$ws.Range("AX49").Value = 0.88129254784759148
$ws.Range("AY49").Value = 0.69816710735424015
$ws.Range("I50").Value = 0.78631925833833027
$ws.Range("AV50").Value = 0.88298208592283345
$ws.Range("AX51").Value = 0.74380949351635894
$ws.Range("AY52").Value = 0.89132554210694093
$ws.Range("AY53").Value = 0.81921456610844468
$ws.Range("AZ53").Value = 0.95129131523099597
$ws.Range("BI53").Value = 0.68528267279987198
$ws.Range("BK53").Value = 0.7103320527783521
$ws.Range("BA55").Value = 0.67534454210030082
$ws.Range("BD55").Value = 0.67405924958275865
$ws.Range("BM55").Value = 0.71020515152697294
$ws.Range("BB56").Value = 0.79922421255856868
$ws.Range("BD58").Value = 0.81148625303188948
$ws.Range("BE58").Value = 0.94346577680085131
$ws.Range("J59").Value = 0.90980101851872841
$ws.Range("AE59").Value = 0.8595823296707048
$ws.Range("BE59").Value = 0.83955954693335166
$ws.Range("BF59").Value = 0.61366347862253234
$ws.Range("BI59").Value = 0.84311233225607862
$ws.Range("BG60").Value = 0.79141931647870678
$ws.Range("BI60").Value = 0.6549088155512176
$ws.Range("BJ60").Value = 0.54156046553107462
$ws.Range("D61").Value = 0.92496377326343304
$ws.Range("BJ61").Value = 0.87147026854400744
$ws.Range("E62").Value = 0.69966863668104207
$ws.Range("BK62").Value = 0.98038380689163362
$ws.Range("AO63").Value = 0.99757659026171308
$ws.Range("BF63").Value = 0.96170450131935792
$ws.Range("BM63").Value = 0.77655443856892248
$ws.Range("BJ64").Value = 0.93218239119941093
$ws.Range("BM64").Value = 0.65861285038026884
$ws.Range("BB66").Value = 0.67023194268471586
$ws.Range("BL66").Value = 0.76183629250886731
$ws.Range("BM66").Value = 0.90930585600269076
$ws.Range("BP66").Value = 0.81203820489915368
$ws.Range("A67").Value = 0.85245237167730203
$ws.Range("B68").Value = 0.79298238981774349
$ws.Range("Z68").Value = 0.82633449145229387
